$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

# The placeholder currently holds a single empty (level-2) bullet paragraph.
# Replace it with two runs so the misspelled-looking "Yange" keeps its own
# run (as PowerPoint does when it flags a word during spell-check) while the
# rest of the sentence is a separate run.
$tr.Text = "Yange"
$tr.LanguageID = "en-US"
[void]$tr.InsertAfter(" Street in Toronto is the pick to open an Italian Restaurant")
